{"js": "// Applies the Abrams_Cover_Letter.docx edits:\n//  1. \"February 20th, 2021\" -> \"February 21st, 2021\" (keeps superscript run)\n//  2. Opening line: \"May Mobility recently posted on their jobs board\"\n//     -> \"GitHub recently posted on LinkedIn\"\n//  3. \"Software Engineer role\" -> \"Software Engineer - Platform role\"\n//  4. \"I am also passionate about creating a greener world\"\n//     -> \"I am passionate about designing maintainable, scalable solutions\"\n//  5. Second mention \"May Mobility\" -> \"GitHub\"\n//  6. \"A sustainable future is a cause I am aligned with, and I would thrive\n//      if given the opportunity to help achieve such a goal\"\n//     -> \"Your collaborative, developer-oriented culture is one I could see\n//         myself thrive in\"\n//  7. \"success of your autonomous mobility solutions\"\n//     -> \"successful operation of your version control and code management\n//         services\"  (keeps the _GoBack bookmark sitting right before the\n//         final period)\n\nconst body = context.document.body;\n\nasync function replaceOnce(searchText, newText, options) {\n  const searchOptions = Object.assign({ matchCase: true }, options || {});\n  const results = body.search(searchText, searchOptions);\n  results.load(\"text\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"edit.js: search text not found: \" + searchText);\n  }\n\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 1a. Day-of-month number.\nawait replaceOnce(\"February 20\", \"February 21\");\n\n// 1b. The superscript ordinal suffix (\"th\" -> \"st\") lives only in the first\n// paragraph, so scope the search there to avoid the many other \"th\"\n// substrings later in the letter (e.g. \"their\", \"Technology\").\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\nconst dateParagraph = paragraphs.items[0];\nconst suffixResults = dateParagraph.search(\"th\", { matchCase: true });\nsuffixResults.load(\"text\");\nawait context.sync();\nif (suffixResults.items.length === 0) {\n  throw new Error(\"edit.js: ordinal suffix 'th' not found in date line\");\n}\nsuffixResults.items[0].insertText(\"st\", Word.InsertLocation.replace);\nawait context.sync();\n\n// 2. Opening sentence: swap the employer + where-it-was-posted phrase.\nawait replaceOnce(\n  \"May Mobility recently posted on their jobs board\",\n  \"GitHub recently posted on LinkedIn\"\n);\n\n// 3. Add the \" - Platform\" qualifier to the job title.\nawait replaceOnce(\"Software Engineer role\", \"Software Engineer - Platform role\");\n\n// 4. Swap the \"passionate about\" clause.\nawait replaceOnce(\n  \"I am also passionate about creating a greener world\",\n  \"I am passionate about designing maintainable, scalable solutions\"\n);\n\n// 5. Second standalone mention of the company name (leave the trailing\n// \". \" run untouched, matching the source formatting).\nawait replaceOnce(\"May Mobility\", \"GitHub\");\n\n// 6. Replace the \"sustainable future\" sentence with the culture sentence.\nawait replaceOnce(\n  \"A sustainable future is a cause I am aligned with, and I would thrive if given the opportunity to help achieve such a goal\",\n  \"Your collaborative, developer-oriented culture is one I could see myself thrive in\"\n);\n\n// 7. Closing paragraph: swap what the company is recognized for.\nawait replaceOnce(\n  \"success of your autonomous mobility solutions\",\n  \"successful operation of your version control and code management services\"\n);\n", "ps1": "# Applies the Abrams_Cover_Letter.docx edits:\n#  1. \"February 20th, 2021\" -> \"February 21st, 2021\" (keeps superscript run)\n#  2. Opening line: \"May Mobility recently posted on their jobs board\"\n#     -> \"GitHub recently posted on LinkedIn\"\n#  3. \"Software Engineer role\" -> \"Software Engineer - Platform role\"\n#  4. \"I am also passionate about creating a greener world\"\n#     -> \"I am passionate about designing maintainable, scalable solutions\"\n#  5. Second mention \"May Mobility\" -> \"GitHub\"\n#  6. \"A sustainable future is a cause I am aligned with, and I would thrive\n#      if given the opportunity to help achieve such a goal\"\n#     -> \"Your collaborative, developer-oriented culture is one I could see\n#         myself thrive in\"\n#  7. \"success of your autonomous mobility solutions\"\n#     -> \"successful operation of your version control and code management\n#         services\"  (keeps the _GoBack bookmark sitting right before the\n#         final period)\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text($findText, $replaceText) {\n    # wdFindContinue = 1, wdReplaceOne = 1 (last arg) -> only the first hit\n    # in the supplied range is touched, just like a single manual Find &\n    # Replace \"Replace\" click.\n    $found = $d.Content.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 1)\n    if (-not $found) {\n        throw \"edit.ps1: text not found: $findText\"\n    }\n}\n\n# 1a. Day-of-month number.\nReplace-Text \"February 20\" \"February 21\"\n\n# 1b. Ordinal suffix \"th\" -> \"st\". Scope the Find to the first paragraph's\n# range only, since \"th\" also occurs as a substring elsewhere in the\n# letter (e.g. \"their\", \"Technology\") and we must only touch the\n# superscripted ordinal suffix on the date line.\n$dateParagraph = $d.Paragraphs.Item(1).Range\n$found = $dateParagraph.Find.Execute(\"th\", $false, $false, $false, $false, $false, $true, 1, $false, \"st\", 1)\nif (-not $found) {\n    throw \"edit.ps1: ordinal suffix 'th' not found in date line\"\n}\n\n# 2. Opening sentence: swap the employer + where-it-was-posted phrase.\nReplace-Text \"May Mobility recently posted on their jobs board\" \"GitHub recently posted on LinkedIn\"\n\n# 3. Add the \" - Platform\" qualifier to the job title.\nReplace-Text \"Software Engineer role\" \"Software Engineer - Platform role\"\n\n# 4. Swap the \"passionate about\" clause.\nReplace-Text \"I am also passionate about creating a greener world\" \"I am passionate about designing maintainable, scalable solutions\"\n\n# 5. Second standalone mention of the company name (leave the trailing\n# \". \" text untouched, matching the source formatting).\nReplace-Text \"May Mobility\" \"GitHub\"\n\n# 6. Replace the \"sustainable future\" sentence with the culture sentence.\nReplace-Text \"A sustainable future is a cause I am aligned with, and I would thrive if given the opportunity to help achieve such a goal\" \"Your collaborative, developer-oriented culture is one I could see myself thrive in\"\n\n# 7. Closing paragraph: swap what the company is recognized for.\nReplace-Text \"success of your autonomous mobility solutions\" \"successful operation of your version control and code management services\"\n"}
